$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 109.94574705276446
$ws.Range("C2").Value = 40.736074395676106
$ws.Range("D2").Value = 48.87056625032141
$ws.Range("E2").Value = 38.110677903548016

$ws.Range("B3").Value = 68.291712523665268
$ws.Range("C3").Value = 44.34243611499398
$ws.Range("D3").Value = 47.032783547921092
$ws.Range("E3").Value = 17.959768971462918

$ws.Range("B1:E3").Select()
